$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Dlk1"
$ws.Range("C2").Value2 = "Notch2"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.028147
$ws.Range("H2").Value2 = 0.084441
$ws.Range("I2").Value2 = 0.0007347845853173872
$ws.Range("J2").Value2 = 0.000734784585317387
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 2.993142333333334
$ws.Range("N2").Value2 = 8.979427000000001
$ws.Range("O2").Value2 = 0.03484385887642424
$ws.Range("P2").Value2 = 0.03484385887642424
$ws.Range("Q2").Value2 = 0.08424797725633335
$ws.Range("R2").Value2 = 0.7582317953070001
$ws.Range("S2").Value2 = 0.00002560273039537094
$ws.Range("T2").Value2 = 0.00002560273039537094

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Dlk1"
$ws.Range("C3").Value2 = "Notch2"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.028147
$ws.Range("H3").Value2 = 0.084441
$ws.Range("I3").Value2 = 0.0007347845853173872
$ws.Range("J3").Value2 = 0.000734784585317387
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 31.995262
$ws.Range("N3").Value2 = 95.985786
$ws.Range("O3").Value2 = 0.3724642097459734
$ws.Range("P3").Value2 = 0.3724642097459735
$ws.Range("Q3").Value2 = 0.9005706395140001
$ws.Range("R3").Value2 = 8.105135755626
$ws.Range("S3").Value2 = 0.0002736809599037634
$ws.Range("T3").Value2 = 0.0002736809599037634

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Dlk1"
$ws.Range("C4").Value2 = "Notch2"
$ws.Range("D4").Value2 = "M2"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.028147
$ws.Range("H4").Value2 = 0.084441
$ws.Range("I4").Value2 = 0.0007347845853173872
$ws.Range("J4").Value2 = 0.000734784585317387
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 34.28929533333334
$ws.Range("N4").Value2 = 102.867886
$ws.Range("O4").Value2 = 0.3991695798295478
$ws.Range("P4").Value2 = 0.3991695798295478
$ws.Range("Q4").Value2 = 0.9651407957473335
$ws.Range("R4").Value2 = 8.686267161726
$ws.Range("S4").Value2 = 0.0002933036541863699
$ws.Range("T4").Value2 = 0.0002933036541863699

$ws.Range("A5").Value2 = "ECs"
$ws.Range("B5").Value2 = "Dlk1"
$ws.Range("C5").Value2 = "Notch2"
$ws.Range("D5").Value2 = "sCs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.028147
$ws.Range("H5").Value2 = 0.084441
$ws.Range("I5").Value2 = 0.0007347845853173872
$ws.Range("J5").Value2 = 0.000734784585317387
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 16.62387466666667
$ws.Range("N5").Value2 = 49.871624
$ws.Range("O5").Value2 = 0.1935223515480544
$ws.Range("P5").Value2 = 0.1935223515480545
$ws.Range("Q5").Value2 = 0.4679122002426667
$ws.Range("R5").Value2 = 4.211209802183999
$ws.Range("S5").Value2 = 0.0001421972408318828
$ws.Range("T5").Value2 = 0.0001421972408318828

$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Dlk1"
$ws.Range("C6").Value2 = "Notch2"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 36.97491766666666
$ws.Range("H6").Value2 = 110.924753
$ws.Range("I6").Value2 = 0.9652396186039789
$ws.Range("J6").Value2 = 0.9652396186039789
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.993142333333334
$ws.Range("N6").Value2 = 8.979427000000001
$ws.Range("O6").Value2 = 0.03484385887642424
$ws.Range("P6").Value2 = 0.03484385887642424
$ws.Range("Q6").Value2 = 110.6711913396146
$ws.Range("R6").Value2 = 996.0407220565311
$ws.Range("S6").Value2 = 0.0336326730525706
$ws.Range("T6").Value2 = 0.0336326730525706

$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Dlk1"
$ws.Range("C7").Value2 = "Notch2"
$ws.Range("D7").Value2 = "FAPs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 36.97491766666666
$ws.Range("H7").Value2 = 110.924753
$ws.Range("I7").Value2 = 0.9652396186039789
$ws.Range("J7").Value2 = 0.9652396186039789
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 31.995262
$ws.Range("N7").Value2 = 95.985786
$ws.Range("O7").Value2 = 0.3724642097459734
$ws.Range("P7").Value2 = 0.3724642097459735
$ws.Range("Q7").Value2 = 1183.022178173429
$ws.Range("R7").Value2 = 10647.19960356086
$ws.Range("S7").Value2 = 0.3595172117588358
$ws.Range("T7").Value2 = 0.3595172117588358

$ws.Range("A8").Value2 = "FAPs"
$ws.Range("B8").Value2 = "Dlk1"
$ws.Range("C8").Value2 = "Notch2"
$ws.Range("D8").Value2 = "M2"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 36.97491766666666
$ws.Range("H8").Value2 = 110.924753
$ws.Range("I8").Value2 = 0.9652396186039789
$ws.Range("J8").Value2 = 0.9652396186039789
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 34.28929533333334
$ws.Range("N8").Value2 = 102.867886
$ws.Range("O8").Value2 = 0.3991695798295478
$ws.Range("P8").Value2 = 0.3991695798295478
$ws.Range("Q8").Value2 = 1267.843871798018
$ws.Range("R8").Value2 = 11410.59484618216
$ws.Range("S8").Value2 = 0.3852942929929832
$ws.Range("T8").Value2 = 0.3852942929929832

$ws.Range("A9").Value2 = "FAPs"
$ws.Range("B9").Value2 = "Dlk1"
$ws.Range("C9").Value2 = "Notch2"
$ws.Range("D9").Value2 = "sCs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 36.97491766666666
$ws.Range("H9").Value2 = 110.924753
$ws.Range("I9").Value2 = 0.9652396186039789
$ws.Range("J9").Value2 = 0.9652396186039789
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 16.62387466666667
$ws.Range("N9").Value2 = 49.871624
$ws.Range("O9").Value2 = 0.1935223515480544
$ws.Range("P9").Value2 = 0.1935223515480545
$ws.Range("Q9").Value2 = 614.6663971009857
$ws.Range("R9").Value2 = 5531.997573908871
$ws.Range("S9").Value2 = 0.1867954407995892
$ws.Range("T9").Value2 = 0.1867954407995892

$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Dlk1"
$ws.Range("C10").Value2 = "Notch2"
$ws.Range("D10").Value2 = "ECs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 1.303400333333333
$ws.Range("H10").Value2 = 3.910201
$ws.Range("I10").Value2 = 0.03402559681070371
$ws.Range("J10").Value2 = 0.03402559681070371
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 2.993142333333334
$ws.Range("N10").Value2 = 8.979427000000001
$ws.Range("O10").Value2 = 0.03484385887642424
$ws.Range("P10").Value2 = 0.03484385887642424
$ws.Range("Q10").Value2 = 3.901262714980778
$ws.Range("R10").Value2 = 35.111364434827
$ws.Range("S10").Value2 = 0.001185583093458271
$ws.Range("T10").Value2 = 0.001185583093458271

$ws.Range("A11").Value2 = "sCs"
$ws.Range("B11").Value2 = "Dlk1"
$ws.Range("C11").Value2 = "Notch2"
$ws.Range("D11").Value2 = "FAPs"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 1.303400333333333
$ws.Range("H11").Value2 = 3.910201
$ws.Range("I11").Value2 = 0.03402559681070371
$ws.Range("J11").Value2 = 0.03402559681070371
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 31.995262
$ws.Range("N11").Value2 = 95.985786
$ws.Range("O11").Value2 = 0.3724642097459734
$ws.Range("P11").Value2 = 0.3724642097459735
$ws.Range("Q11").Value2 = 41.70263515588733
$ws.Range("R11").Value2 = 375.323716402986
$ws.Range("S11").Value2 = 0.01267331702723387
$ws.Range("T11").Value2 = 0.01267331702723387

$ws.Range("A12").Value2 = "sCs"
$ws.Range("B12").Value2 = "Dlk1"
$ws.Range("C12").Value2 = "Notch2"
$ws.Range("D12").Value2 = "M2"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 1.303400333333333
$ws.Range("H12").Value2 = 3.910201
$ws.Range("I12").Value2 = 0.03402559681070371
$ws.Range("J12").Value2 = 0.03402559681070371
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 34.28929533333334
$ws.Range("N12").Value2 = 102.867886
$ws.Range("O12").Value2 = 0.3991695798295478
$ws.Range("P12").Value2 = 0.3991695798295478
$ws.Range("Q12").Value2 = 44.69267896723178
$ws.Range("R12").Value2 = 402.234110705086
$ws.Range("S12").Value2 = 0.0135819831823782
$ws.Range("T12").Value2 = 0.0135819831823782

$ws.Range("A13").Value2 = "sCs"
$ws.Range("B13").Value2 = "Dlk1"
$ws.Range("C13").Value2 = "Notch2"
$ws.Range("D13").Value2 = "sCs"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 1.303400333333333
$ws.Range("H13").Value2 = 3.910201
$ws.Range("I13").Value2 = 0.03402559681070371
$ws.Range("J13").Value2 = 0.03402559681070371
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 16.62387466666667
$ws.Range("N13").Value2 = 49.871624
$ws.Range("O13").Value2 = 0.1935223515480544
$ws.Range("P13").Value2 = 0.1935223515480545
$ws.Range("Q13").Value2 = 21.66756378182489
$ws.Range("R13").Value2 = 195.008074036424
$ws.Range("S13").Value2 = 0.006584713507633364
$ws.Range("T13").Value2 = 0.006584713507633365

